$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "72.986.79"
$ws.Range("E2").Value = "  +3.14%  "

Set-TextValue $ws.Range("D3") "3.986.79"
$ws.Range("E3").Value = "  +1.23%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue $ws.Range("D5") "596.96"
$ws.Range("E5").Value = "  +11.33%  "

Set-TextValue $ws.Range("D6") "159.89"
$ws.Range("E6").Value = "  +8.06%  "

Set-TextValue $ws.Range("D7") "0.683"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  -0.03%  "

Set-TextValue $ws.Range("D9") "0.749"
$ws.Range("E9").Value = "  +1.70%  "

Set-TextValue $ws.Range("D10") "0.169"
$ws.Range("E10").Value = "  +2.44%  "

Set-TextValue $ws.Range("D11") "53.84"
$ws.Range("E11").Value = "  -3.10%  "

Set-TextValue $ws.Range("D12") "0.0000319"
$ws.Range("E12").Value = "  +1.59%  "

Set-TextValue $ws.Range("D13") "10.96"
$ws.Range("E13").Value = "  +3.52%  "

Set-TextValue $ws.Range("D14") "4.635.63"
$ws.Range("E14").Value = "  +0.86%  "

Set-TextValue $ws.Range("D15") "3.995.91"
$ws.Range("E15").Value = "  +0.97%  "

Set-TextValue $ws.Range("D16") "1.26"
$ws.Range("E16").Value = "  +8.55%  "

Set-TextValue $ws.Range("D17") "14.10"
$ws.Range("E17").Value = "  +2.36%  "

Set-TextValue $ws.Range("D18") "20.32"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("E19").Value = "  +0.35%  "

Set-TextValue $ws.Range("D20") "72.764.35"
$ws.Range("E20").Value = "  +2.84%  "

Set-TextValue $ws.Range("D21") "434.27"
$ws.Range("E21").Value = "  +2.39%  "

Set-TextValue $ws.Range("D22") "4.78"
$ws.Range("E22").Value = "  +13.73%  "

Set-TextValue $ws.Range("D23") "96.12"
$ws.Range("E23").Value = "  -0.69%  "

Set-TextValue $ws.Range("D24") "3.42"
$ws.Range("E24").Value = "  -4.03%  "

Set-TextValue $ws.Range("D25") "14.20"
$ws.Range("E25").Value = "  -1.16%  "

Set-TextValue $ws.Range("D26") "4.38"
$ws.Range("E26").Value = "  +16.77%  "

Set-TextValue $ws.Range("D27") "11.21"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("E28").Value = "  +1.07%  "

Set-TextValue $ws.Range("D29") "10.46"
$ws.Range("E29").Value = "  -1.32%  "

Set-TextValue $ws.Range("D30") "36.33"
$ws.Range("E30").Value = "  +0.07%  "

Set-TextValue $ws.Range("D31") "7.80"
$ws.Range("E31").Value = "  +0.80%  "

Set-TextValue $ws.Range("D32") "13.76"
$ws.Range("E32").Value = "  +3.45%  "

Set-TextValue $ws.Range("D33") "0.130"
$ws.Range("E33").Value = "  +0.09%  "

Set-TextValue $ws.Range("D34") "48.06"
$ws.Range("E34").Value = "  -3.73%  "

Set-TextValue $ws.Range("D35") "669.08"
$ws.Range("E35").Value = "  -1.94%  "

Set-TextValue $ws.Range("D36") "70.75"
$ws.Range("E36").Value = "  +8.79%  "

Set-TextValue $ws.Range("D37") "0.0₃0903"
$ws.Range("E37").Value = "  +10.91%  "

Set-TextValue $ws.Range("D38") "0.436"
$ws.Range("E38").Value = "  +0.43%  "

$ws.Range("E39").Value = "  -0.11%  "

Set-TextValue $ws.Range("D40") "0.145"
$ws.Range("E40").Value = "  -2.41%  "

Set-TextValue $ws.Range("D41") "3.34"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  +4.43%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D44") "0.0493"
$ws.Range("E44").Value = "  +2.90%  "

Set-TextValue $ws.Range("D45") "10.59"
$ws.Range("E45").Value = "  +8.53%  "

Set-TextValue $ws.Range("D46") "0.149"
$ws.Range("E46").Value = "  +0.89%  "

Set-TextValue $ws.Range("D47") "3.42"
$ws.Range("E47").Value = "  +2.92%  "

Set-TextValue $ws.Range("D48") "2.60"
$ws.Range("E48").Value = "  -2.76%  "

Set-TextValue $ws.Range("D49") "2.870.73"
$ws.Range("E49").Value = "  +8.82%  "

Set-TextValue $ws.Range("D50") "3.03"
$ws.Range("E50").Value = "  +1.70%  "

Set-TextValue $ws.Range("D51") "3.38"
$ws.Range("E51").Value = "  +4.53%  "
